$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Updated "Estado de Cuenta" (account statement) data table.
# Column layout: B=Tipo Doc, C=N Doc Trabajador, D=Nombre Trabajador,
#                E=Periodo Mora, F=Valor Mora, G=Salario Basico
# The workers' overdue periods are reorganized/expanded (chronological order,
# part 1 of new statement), replacing the previous 5-period block for a
# single worker with an interleaved two-worker table.

$rows = @(
    @{ Row = 16; Tipo = "CC"; Doc = "86050699";   Nombre = "CARLOS ENRIQUE GODOY RIAÑO";   Periodo = "1810"; Mora = 31249; Salario = 781242 },
    @{ Row = 17; Tipo = "CC"; Doc = "86050699";   Nombre = "CARLOS ENRIQUE GODOY RIAÑO";   Periodo = "1811"; Mora = 31249; Salario = 781242 },
    @{ Row = 18; Tipo = "CC"; Doc = "1143360772"; Nombre = "ANGELICA SUSANA GARCIA PETRO"; Periodo = "1811"; Mora = 31249; Salario = 781242 },
    @{ Row = 19; Tipo = "CC"; Doc = "86050699";   Nombre = "CARLOS ENRIQUE GODOY RIAÑO";   Periodo = "1901"; Mora = 31249; Salario = 781242 },
    @{ Row = 20; Tipo = "CC"; Doc = "1143360772"; Nombre = "ANGELICA SUSANA GARCIA PETRO"; Periodo = "1901"; Mora = 31249; Salario = 781242 },
    @{ Row = 21; Tipo = "CC"; Doc = "86050699";   Nombre = "CARLOS ENRIQUE GODOY RIAÑO";   Periodo = "1902"; Mora = 31249; Salario = 781242 },
    @{ Row = 22; Tipo = "CC"; Doc = "1143360772"; Nombre = "ANGELICA SUSANA GARCIA PETRO"; Periodo = "1902"; Mora = 31249; Salario = 781242 },
    @{ Row = 23; Tipo = "CC"; Doc = "86050699";   Nombre = "CARLOS ENRIQUE GODOY RIAÑO";   Periodo = "1903"; Mora = 26041; Salario = 781242 },
    @{ Row = 24; Tipo = "CC"; Doc = "1143360772"; Nombre = "ANGELICA SUSANA GARCIA PETRO"; Periodo = "1903"; Mora = 26041; Salario = 781242 }
)

foreach ($r in $rows) {
    $n = $r.Row
    $ws.Cells.Item($n, 2).Value = $r.Tipo
    $ws.Cells.Item($n, 3).Value = $r.Doc
    $ws.Cells.Item($n, 4).Value = $r.Nombre
    $ws.Cells.Item($n, 5).Value = $r.Periodo
    $ws.Cells.Item($n, 6).Value = $r.Mora
    $ws.Cells.Item($n, 7).Value = $r.Salario
}
